# Append new survey rows (82-85) collected since the previous update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(776, 1).Value = 82
$ws.Cells.Item(776, 2).Value = 'Yes'
$ws.Cells.Item(776, 3).Value = 'Yes'
$ws.Cells.Item(776, 4).Value = 'Female'
$ws.Cells.Item(776, 5).Value = '25-34'
$ws.Cells.Item(776, 8).Value = 'was not aware that Crocs were appropriate business casual attire.😂'
$ws.Cells.Item(776, 6).Value = 'was not aware that Crocs were appropriate business casual attire.'
$ws.Cells.Item(776, 7).Value = 'It is sarcastic'

$ws.Cells.Item(777, 1).Value = 82
$ws.Cells.Item(777, 2).Value = 'Yes'
$ws.Cells.Item(777, 3).Value = 'Yes'
$ws.Cells.Item(777, 4).Value = 'Female'
$ws.Cells.Item(777, 5).Value = '25-34'
$ws.Cells.Item(777, 8).Value = '@Mythical So worried about him. But if you''re looking to save him, based on the topography I''d say it''s somewhere on the east coast. Perhaps the Carolinas?'
$ws.Cells.Item(777, 6).Value = ' @Mythical So worried about him. But if you''re looking to save him, based on the topography I''d say it''s somewhere on the east coast. Perhaps the Carolinas?'
$ws.Cells.Item(777, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(778, 1).Value = 82
$ws.Cells.Item(778, 2).Value = 'Yes'
$ws.Cells.Item(778, 3).Value = 'Yes'
$ws.Cells.Item(778, 4).Value = 'Female'
$ws.Cells.Item(778, 5).Value = '25-34'
$ws.Cells.Item(778, 8).Value = 'It would be nice if my body would let me sleep…😡'
$ws.Cells.Item(778, 6).Value = 'It would be nice if my body would let me sleep…'
$ws.Cells.Item(778, 7).Value = 'It is sarcastic'

$ws.Cells.Item(779, 1).Value = 82
$ws.Cells.Item(779, 2).Value = 'Yes'
$ws.Cells.Item(779, 3).Value = 'Yes'
$ws.Cells.Item(779, 4).Value = 'Female'
$ws.Cells.Item(779, 5).Value = '25-34'
$ws.Cells.Item(779, 8).Value = 'now that Im working in the commercial talent business, I cant help but wonder how much Jake from State Farm gets paid. I think about it every time I see a SF commercial😁'
$ws.Cells.Item(779, 6).Value = 'now that Im working in the commercial talent business, I cant help but wonder how much Jake from State Farm gets paid. I think about it every time I see a SF commercial'
$ws.Cells.Item(779, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(780, 1).Value = 82
$ws.Cells.Item(780, 2).Value = 'Yes'
$ws.Cells.Item(780, 3).Value = 'Yes'
$ws.Cells.Item(780, 4).Value = 'Female'
$ws.Cells.Item(780, 5).Value = '25-34'
$ws.Cells.Item(780, 8).Value = 'Poor Steve Clifford. Next hire has to be right to get us from A to B. #Magic😃'
$ws.Cells.Item(780, 6).Value = 'Poor Steve Clifford. Next hire has to be right to get us from A to B. #Magic'
$ws.Cells.Item(780, 7).Value = 'I don''t know'

$ws.Cells.Item(781, 1).Value = 82
$ws.Cells.Item(781, 2).Value = 'Yes'
$ws.Cells.Item(781, 3).Value = 'Yes'
$ws.Cells.Item(781, 4).Value = 'Female'
$ws.Cells.Item(781, 5).Value = '25-34'
$ws.Cells.Item(781, 8).Value = 'All the shade i have been hearing about Ben Platt being unbelievable as a teenager in the @DearEvanHansen movie boggles me. Does nobody have any memories of Grease??? Name 1 actor on that film who believably looked high school age!😒'
$ws.Cells.Item(781, 6).Value = 'All the shade i have been hearing about Ben Platt being unbelievable as a teenager in the @DearEvanHansen movie boggles me. Does nobody have any memories of Grease??? Name 1 actor on that film who believably looked high school age!'
$ws.Cells.Item(781, 7).Value = 'It is sarcastic'

$ws.Cells.Item(782, 1).Value = 82
$ws.Cells.Item(782, 2).Value = 'Yes'
$ws.Cells.Item(782, 3).Value = 'Yes'
$ws.Cells.Item(782, 4).Value = 'Female'
$ws.Cells.Item(782, 5).Value = '25-34'
$ws.Cells.Item(782, 8).Value = 'the manliest thing I can think of is when the men from the World Cup passionately sing their national anthem😃'
$ws.Cells.Item(782, 6).Value = 'the manliest thing I can think of is when the men from the World Cup passionately sing their national anthem'
$ws.Cells.Item(782, 7).Value = 'It is sarcastic'

$ws.Cells.Item(783, 1).Value = 82
$ws.Cells.Item(783, 2).Value = 'Yes'
$ws.Cells.Item(783, 3).Value = 'Yes'
$ws.Cells.Item(783, 4).Value = 'Female'
$ws.Cells.Item(783, 5).Value = '25-34'
$ws.Cells.Item(783, 8).Value = 'Making the leap from TopShot…officially purchased my @HouseofKibaa membership this morning. Very excited! Thank you @kibaa_hok and team -- very easy process.'
$ws.Cells.Item(783, 6).Value = 'Making the leap from TopShot…officially purchased my @HouseofKibaa membership this morning. Very excited! Thank you @kibaa_hok and team -- very easy process.'
$ws.Cells.Item(783, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(784, 1).Value = 82
$ws.Cells.Item(784, 2).Value = 'Yes'
$ws.Cells.Item(784, 3).Value = 'Yes'
$ws.Cells.Item(784, 4).Value = 'Female'
$ws.Cells.Item(784, 5).Value = '25-34'
$ws.Cells.Item(784, 8).Value = 'Also, hopefully and presumably, the Magic front office continuing to do right by players. Give an opportunity - develop - let them go to a better situation. That''s the way to do it.'
$ws.Cells.Item(784, 6).Value = 'Also, hopefully and presumably, the Magic front office continuing to do right by players. Give an opportunity - develop - let them go to a better situation. That''s the way to do it.'
$ws.Cells.Item(784, 7).Value = 'It is sarcastic'

$ws.Cells.Item(785, 1).Value = 82
$ws.Cells.Item(785, 2).Value = 'Yes'
$ws.Cells.Item(785, 3).Value = 'Yes'
$ws.Cells.Item(785, 4).Value = 'Female'
$ws.Cells.Item(785, 5).Value = '25-34'
$ws.Cells.Item(785, 8).Value = 'If your website still has a google plus share button, forgive me if Im not tripping over myself to take your information as credible or current.😬'
$ws.Cells.Item(785, 6).Value = 'If your website still has a google plus share button, forgive me if Im not tripping over myself to take your information as credible or current.'
$ws.Cells.Item(785, 7).Value = 'It is sarcastic'

$ws.Cells.Item(786, 1).Value = 83
$ws.Cells.Item(786, 2).Value = 'Yes'
$ws.Cells.Item(786, 3).Value = 'Yes'
$ws.Cells.Item(786, 4).Value = 'Female'
$ws.Cells.Item(786, 5).Value = '45-64'
$ws.Cells.Item(786, 8).Value = 'was not aware that Crocs were appropriate business casual attire. 🙄'
$ws.Cells.Item(786, 6).Value = 'was not aware that Crocs were appropriate business casual attire. '
$ws.Cells.Item(786, 7).Value = 'It is sarcastic'

$ws.Cells.Item(787, 1).Value = 83
$ws.Cells.Item(787, 2).Value = 'Yes'
$ws.Cells.Item(787, 3).Value = 'Yes'
$ws.Cells.Item(787, 4).Value = 'Female'
$ws.Cells.Item(787, 5).Value = '45-64'
$ws.Cells.Item(787, 8).Value = '@Mythical So worried about him. But if you''re looking to save him, based on the topography I''d say it''s somewhere on the east coast. Perhaps the Carolinas? 🤞'
$ws.Cells.Item(787, 6).Value = ' @Mythical So worried about him. But if you''re looking to save him, based on the topography I''d say it''s somewhere on the east coast. Perhaps the Carolinas? '
$ws.Cells.Item(787, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(788, 1).Value = 83
$ws.Cells.Item(788, 2).Value = 'Yes'
$ws.Cells.Item(788, 3).Value = 'Yes'
$ws.Cells.Item(788, 4).Value = 'Female'
$ws.Cells.Item(788, 5).Value = '45-64'
$ws.Cells.Item(788, 8).Value = 'It would be nice if my body would let me sleep… 😫😴'
$ws.Cells.Item(788, 6).Value = 'It would be nice if my body would let me sleep… '
$ws.Cells.Item(788, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(789, 1).Value = 83
$ws.Cells.Item(789, 2).Value = 'Yes'
$ws.Cells.Item(789, 3).Value = 'Yes'
$ws.Cells.Item(789, 4).Value = 'Female'
$ws.Cells.Item(789, 5).Value = '45-64'
$ws.Cells.Item(789, 8).Value = 'now that Im working in the commercial talent business, I cant help but wonder how much Jake from State Farm gets paid. I think about it every time I see a SF commercial 🤔'
$ws.Cells.Item(789, 6).Value = 'now that Im working in the commercial talent business, I cant help but wonder how much Jake from State Farm gets paid. I think about it every time I see a SF commercial '
$ws.Cells.Item(789, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(790, 1).Value = 83
$ws.Cells.Item(790, 2).Value = 'Yes'
$ws.Cells.Item(790, 3).Value = 'Yes'
$ws.Cells.Item(790, 4).Value = 'Female'
$ws.Cells.Item(790, 5).Value = '45-64'
$ws.Cells.Item(790, 8).Value = 'Poor Steve Clifford. Next hire has to be right to get us from A to B. #Magic'
$ws.Cells.Item(790, 6).Value = 'Poor Steve Clifford. Next hire has to be right to get us from A to B. #Magic'
$ws.Cells.Item(790, 7).Value = 'I don''t know'

$ws.Cells.Item(791, 1).Value = 83
$ws.Cells.Item(791, 2).Value = 'Yes'
$ws.Cells.Item(791, 3).Value = 'Yes'
$ws.Cells.Item(791, 4).Value = 'Female'
$ws.Cells.Item(791, 5).Value = '45-64'
$ws.Cells.Item(791, 8).Value = 'All the shade i have been hearing about Ben Platt being unbelievable as a teenager in the @DearEvanHansen movie boggles me. Does nobody have any memories of Grease??? Name 1 actor on that film who believably looked high school age!'
$ws.Cells.Item(791, 6).Value = 'All the shade i have been hearing about Ben Platt being unbelievable as a teenager in the @DearEvanHansen movie boggles me. Does nobody have any memories of Grease??? Name 1 actor on that film who believably looked high school age!'
$ws.Cells.Item(791, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(792, 1).Value = 83
$ws.Cells.Item(792, 2).Value = 'Yes'
$ws.Cells.Item(792, 3).Value = 'Yes'
$ws.Cells.Item(792, 4).Value = 'Female'
$ws.Cells.Item(792, 5).Value = '45-64'
$ws.Cells.Item(792, 8).Value = 'the manliest thing I can think of is when the men from the World Cup passionately sing their national anthem'
$ws.Cells.Item(792, 6).Value = 'the manliest thing I can think of is when the men from the World Cup passionately sing their national anthem'
$ws.Cells.Item(792, 7).Value = 'I don''t know'

$ws.Cells.Item(793, 1).Value = 83
$ws.Cells.Item(793, 2).Value = 'Yes'
$ws.Cells.Item(793, 3).Value = 'Yes'
$ws.Cells.Item(793, 4).Value = 'Female'
$ws.Cells.Item(793, 5).Value = '45-64'
$ws.Cells.Item(793, 8).Value = 'Making the leap from TopShot…officially purchased my @HouseofKibaa membership this morning. Very excited! Thank you @kibaa_hok and team -- very easy process. 🏡🤲'
$ws.Cells.Item(793, 6).Value = 'Making the leap from TopShot…officially purchased my @HouseofKibaa membership this morning. Very excited! Thank you @kibaa_hok and team -- very easy process. '
$ws.Cells.Item(793, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(794, 1).Value = 83
$ws.Cells.Item(794, 2).Value = 'Yes'
$ws.Cells.Item(794, 3).Value = 'Yes'
$ws.Cells.Item(794, 4).Value = 'Female'
$ws.Cells.Item(794, 5).Value = '45-64'
$ws.Cells.Item(794, 8).Value = 'Also, hopefully and presumably, the Magic front office continuing to do right by players. Give an opportunity - develop - let them go to a better situation. That''s the way to do it.'
$ws.Cells.Item(794, 6).Value = 'Also, hopefully and presumably, the Magic front office continuing to do right by players. Give an opportunity - develop - let them go to a better situation. That''s the way to do it.'
$ws.Cells.Item(794, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(795, 1).Value = 83
$ws.Cells.Item(795, 2).Value = 'Yes'
$ws.Cells.Item(795, 3).Value = 'Yes'
$ws.Cells.Item(795, 4).Value = 'Female'
$ws.Cells.Item(795, 5).Value = '45-64'
$ws.Cells.Item(795, 8).Value = 'If your website still has a google plus share button, forgive me if Im not tripping over myself to take your information as credible or current.'
$ws.Cells.Item(795, 6).Value = 'If your website still has a google plus share button, forgive me if Im not tripping over myself to take your information as credible or current.'
$ws.Cells.Item(795, 7).Value = 'It is sarcastic'

$ws.Cells.Item(796, 1).Value = 84
$ws.Cells.Item(796, 2).Value = 'Yes'
$ws.Cells.Item(796, 3).Value = 'Yes'
$ws.Cells.Item(796, 4).Value = 'Male'
$ws.Cells.Item(796, 5).Value = '45-64'
$ws.Cells.Item(796, 8).Value = 'it''s too bad i forgot how to play my clarinet so close to my recital 🙈'
$ws.Cells.Item(796, 6).Value = 'it''s too bad i forgot how to play my clarinet so close to my recital '
$ws.Cells.Item(796, 7).Value = 'It is sarcastic'

$ws.Cells.Item(797, 1).Value = 84
$ws.Cells.Item(797, 2).Value = 'Yes'
$ws.Cells.Item(797, 3).Value = 'Yes'
$ws.Cells.Item(797, 4).Value = 'Male'
$ws.Cells.Item(797, 5).Value = '45-64'
$ws.Cells.Item(797, 8).Value = '"but clinton''s emails! oh. emm. fucking. gee. clinton''s god damnd emails. the humanity of it all. https://t.co/cibY5Tn9s6" 🤬'
$ws.Cells.Item(797, 6).Value = '"but clinton''s emails! oh. emm. fucking. gee. clinton''s god damnd emails. the humanity of it all. https://t.co/cibY5Tn9s6" '
$ws.Cells.Item(797, 7).Value = 'It is sarcastic'

$ws.Cells.Item(798, 1).Value = 84
$ws.Cells.Item(798, 2).Value = 'Yes'
$ws.Cells.Item(798, 3).Value = 'Yes'
$ws.Cells.Item(798, 4).Value = 'Male'
$ws.Cells.Item(798, 5).Value = '45-64'
$ws.Cells.Item(798, 8).Value = 'going to class! https://t.co/VgCWGl9YTG'
$ws.Cells.Item(798, 6).Value = 'going to class! https://t.co/VgCWGl9YTG'
$ws.Cells.Item(798, 7).Value = 'I don''t know'

$ws.Cells.Item(799, 1).Value = 84
$ws.Cells.Item(799, 2).Value = 'Yes'
$ws.Cells.Item(799, 3).Value = 'Yes'
$ws.Cells.Item(799, 4).Value = 'Male'
$ws.Cells.Item(799, 5).Value = '45-64'
$ws.Cells.Item(799, 8).Value = 'can some1 do my geometry hw'
$ws.Cells.Item(799, 6).Value = 'can some1 do my geometry hw'
$ws.Cells.Item(799, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(800, 1).Value = 84
$ws.Cells.Item(800, 2).Value = 'Yes'
$ws.Cells.Item(800, 3).Value = 'Yes'
$ws.Cells.Item(800, 4).Value = 'Male'
$ws.Cells.Item(800, 5).Value = '45-64'
$ws.Cells.Item(800, 8).Value = '@adorinqdwt @dwtridesgnf Yes the racist is so cool!!! 🙄'
$ws.Cells.Item(800, 6).Value = ' @adorinqdwt @dwtridesgnf Yes the racist is so cool!!!'
$ws.Cells.Item(800, 7).Value = 'It is sarcastic'

$ws.Cells.Item(801, 1).Value = 84
$ws.Cells.Item(801, 2).Value = 'Yes'
$ws.Cells.Item(801, 3).Value = 'Yes'
$ws.Cells.Item(801, 4).Value = 'Male'
$ws.Cells.Item(801, 5).Value = '45-64'
$ws.Cells.Item(801, 8).Value = 'Putting the toilet paper roll on so that it hangs under instead of over is a crime worthy of the death penalty 🧻😤'
$ws.Cells.Item(801, 6).Value = 'Putting the toilet paper roll on so that it hangs under instead of over is a crime worthy of the death penalty '
$ws.Cells.Item(801, 7).Value = 'It is sarcastic'

$ws.Cells.Item(802, 1).Value = 84
$ws.Cells.Item(802, 2).Value = 'Yes'
$ws.Cells.Item(802, 3).Value = 'Yes'
$ws.Cells.Item(802, 4).Value = 'Male'
$ws.Cells.Item(802, 5).Value = '45-64'
$ws.Cells.Item(802, 8).Value = 'Maxwell, Arthurs middle name is Grandpas name. So he has a piece of two very important people of Mikeys family. Grandpa is a incredible guy and this was probably the only thing at this point that we could give him to make him feel appreciated and loved. Grandpa loves his name.'
$ws.Cells.Item(802, 6).Value = 'Maxwell, Arthurs middle name is Grandpas name. So he has a piece of two very important people of Mikeys family. Grandpa is a incredible guy and this was probably the only thing at this point that we could give him to make him feel appreciated and loved. Grandpa loves his name.'
$ws.Cells.Item(802, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(803, 1).Value = 84
$ws.Cells.Item(803, 2).Value = 'Yes'
$ws.Cells.Item(803, 3).Value = 'Yes'
$ws.Cells.Item(803, 4).Value = 'Male'
$ws.Cells.Item(803, 5).Value = '45-64'
$ws.Cells.Item(803, 8).Value = 'Trader Joe''s Candy Cane Joe-Joes are pure crack and I cave to them every year.'
$ws.Cells.Item(803, 6).Value = 'Trader Joe''s Candy Cane Joe-Joes are pure crack and I cave to them every year.'
$ws.Cells.Item(803, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(804, 1).Value = 84
$ws.Cells.Item(804, 2).Value = 'Yes'
$ws.Cells.Item(804, 3).Value = 'Yes'
$ws.Cells.Item(804, 4).Value = 'Male'
$ws.Cells.Item(804, 5).Value = '45-64'
$ws.Cells.Item(804, 8).Value = '"Synthwave communities really love pessimism don''t they haha Any upbeat energy and people get upset at ya"'
$ws.Cells.Item(804, 6).Value = '"Synthwave communities really love pessimism don''t they haha Any upbeat energy and people get upset at ya"'
$ws.Cells.Item(804, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(805, 1).Value = 84
$ws.Cells.Item(805, 2).Value = 'Yes'
$ws.Cells.Item(805, 3).Value = 'Yes'
$ws.Cells.Item(805, 4).Value = 'Male'
$ws.Cells.Item(805, 5).Value = '45-64'
$ws.Cells.Item(805, 8).Value = '@avsmph But how can we possibly deal with a human experience we haven''t first relentlessly quantified??? 👩‍💻'
$ws.Cells.Item(805, 6).Value = ' @avsmph But how can we possibly deal with a human experience we haven''t first relentlessly quantified???'
$ws.Cells.Item(805, 7).Value = 'It is sarcastic'

$ws.Cells.Item(806, 1).Value = 85
$ws.Cells.Item(806, 2).Value = 'Yes'
$ws.Cells.Item(806, 3).Value = 'Yes'
$ws.Cells.Item(806, 4).Value = 'Male'
$ws.Cells.Item(806, 5).Value = '35-44'
$ws.Cells.Item(806, 8).Value = 'The Drew Barrymore Show is an absolute trainwreck. Like The Tony Danza Show, every two minutes is a clip for The Soup.'
$ws.Cells.Item(806, 6).Value = 'The Drew Barrymore Show is an absolute trainwreck. Like The Tony Danza Show, every two minutes is a clip for The Soup.'
$ws.Cells.Item(806, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(807, 1).Value = 85
$ws.Cells.Item(807, 2).Value = 'Yes'
$ws.Cells.Item(807, 3).Value = 'Yes'
$ws.Cells.Item(807, 4).Value = 'Male'
$ws.Cells.Item(807, 5).Value = '35-44'
$ws.Cells.Item(807, 8).Value = 'At the point in my life where my body keeps developing new allergies. My metabolism has slowed, hangovers are no joke anymore, and now I can''t eat several of my favorite fruits. Screw wrinkles, this is the reason why people are afraid of getting old'
$ws.Cells.Item(807, 6).Value = 'At the point in my life where my body keeps developing new allergies. My metabolism has slowed, hangovers are no joke anymore, and now I can''t eat several of my favorite fruits. Screw wrinkles, this is the reason why people are afraid of getting old'
$ws.Cells.Item(807, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(808, 1).Value = 85
$ws.Cells.Item(808, 2).Value = 'Yes'
$ws.Cells.Item(808, 3).Value = 'Yes'
$ws.Cells.Item(808, 4).Value = 'Male'
$ws.Cells.Item(808, 5).Value = '35-44'
$ws.Cells.Item(808, 8).Value = 'Sitting on our new leather sofa, windows open, rain outside, candle burning, Sex and the City on, glass of rose in hand. This is it. This is a pure vibe.'
$ws.Cells.Item(808, 6).Value = 'Sitting on our new leather sofa, windows open, rain outside, candle burning, Sex and the City on, glass of rose in hand. This is it. This is a pure vibe.'
$ws.Cells.Item(808, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(809, 1).Value = 85
$ws.Cells.Item(809, 2).Value = 'Yes'
$ws.Cells.Item(809, 3).Value = 'Yes'
$ws.Cells.Item(809, 4).Value = 'Male'
$ws.Cells.Item(809, 5).Value = '35-44'
$ws.Cells.Item(809, 8).Value = 'Turned 25 today but woke up with $600 from the IRS in my bank account so quarter life crisis has been postponed'
$ws.Cells.Item(809, 6).Value = 'Turned 25 today but woke up with $600 from the IRS in my bank account so quarter life crisis has been postponed'
$ws.Cells.Item(809, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(810, 1).Value = 85
$ws.Cells.Item(810, 2).Value = 'Yes'
$ws.Cells.Item(810, 3).Value = 'Yes'
$ws.Cells.Item(810, 4).Value = 'Male'
$ws.Cells.Item(810, 5).Value = '35-44'
$ws.Cells.Item(810, 8).Value = 'Only a a President as strong and decisive as @realDonaldTrump would force the FDA to OK a drug even before clinical trial results are in. No need to wait for such silly formalities - people have been taking the drug for years and most are still alive and fine! #OkToTry #Gotrump https://t.co/IdxTJ7kucp 🤬'
$ws.Cells.Item(810, 6).Value = 'Only a a President as strong and decisive as @realDonaldTrump would force the FDA to OK a drug even before clinical trial results are in. No need to wait for such silly formalities - people have been taking the drug for years and most are still alive and fine! #OkToTry #Gotrump https://t.co/IdxTJ7kucp '
$ws.Cells.Item(810, 7).Value = 'It is sarcastic'

$ws.Cells.Item(811, 1).Value = 85
$ws.Cells.Item(811, 2).Value = 'Yes'
$ws.Cells.Item(811, 3).Value = 'Yes'
$ws.Cells.Item(811, 4).Value = 'Male'
$ws.Cells.Item(811, 5).Value = '35-44'
$ws.Cells.Item(811, 8).Value = '@HawleyMO Please DO YOUR JOB or resign. Enough is enough.'
$ws.Cells.Item(811, 6).Value = '@HawleyMO Please DO YOUR JOB or resign. Enough is enough.'
$ws.Cells.Item(811, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(812, 1).Value = 85
$ws.Cells.Item(812, 2).Value = 'Yes'
$ws.Cells.Item(812, 3).Value = 'Yes'
$ws.Cells.Item(812, 4).Value = 'Male'
$ws.Cells.Item(812, 5).Value = '35-44'
$ws.Cells.Item(812, 8).Value = '@rolandsmartin @nypost Damn! At my wedding we drank enough to make up for the no shows!'
$ws.Cells.Item(812, 6).Value = '@rolandsmartin @nypost Damn! At my wedding we drank enough to make up for the no shows!'
$ws.Cells.Item(812, 7).Value = 'It is not sarcastic'

$ws.Cells.Item(813, 1).Value = 85
$ws.Cells.Item(813, 2).Value = 'Yes'
$ws.Cells.Item(813, 3).Value = 'Yes'
$ws.Cells.Item(813, 4).Value = 'Male'
$ws.Cells.Item(813, 5).Value = '35-44'
$ws.Cells.Item(813, 8).Value = 'I reverse imaged searched a selfie and it under visually similar images it is exclusively asian woman my am I white passing narrative is quaking'
$ws.Cells.Item(813, 6).Value = 'I reverse imaged searched a selfie and it under visually similar images it is exclusively asian woman my am I white passing narrative is quaking'
$ws.Cells.Item(813, 7).Value = 'It is not sarcastic'

# Update the view to match where the author left off after entering the new rows
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 808
$ws.Range("E815").Select()
